$d = $word.ActiveDocument

# Locate the target paragraph robustly: the document contains two "As an
# external company" user-story headers; we only want to fix the second one,
# i.e. the one immediately followed by the "I want to be able to adjust
# prices ..." story. Identify it by that surrounding context rather than a
# hard-coded paragraph index.
$count = $d.Paragraphs.Count
$targetIndex = 0
for ($i = 1; $i -lt $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("As an external company")) {
        $nextText = $d.Paragraphs($i + 1).Range.Text
        if ($nextText.StartsWith("I want to be able to adjust prices")) {
            $targetIndex = $i
            break
        }
    }
}

$para = $d.Paragraphs($targetIndex)
$pr = $para.Range
$s = $pr.Start

# The paragraph reads "As an external company" built from four runs:
#   "As a" | "n " | "external " (bold) | "company" (bold)
# The fix drops the "n" and the whole "external " run, leaving "As a company"
# (single space), i.e. the second run's text becomes " " and the third run
# disappears entirely while the first and last runs stay untouched.

# 1) Delete the bold "external " run outright (positions 6..15 relative to
#    the paragraph start: "As a" + "n " = 6 chars in).
$rExternal = $d.Range($s + 6, $s + 15)
$rExternal.Delete()

# 2) Shrink the "n " run down to just " ". Editing this run's text in place
#    would normally let the engine silently coalesce it into the preceding
#    "As a" run (both share identical non-bold run formatting). Nudge "As a"
#    to a different, temporary format first so the edited run cannot merge
#    into it, then restore "As a" afterwards.
$rAsA = $d.Range($s, $s + 4)
$rAsA.Bold = 1

$rN = $d.Range($s + 4, $s + 6)
$rN.Text = " "

$rAsA2 = $d.Range($s, $s + 4)
$rAsA2.Bold = 0

Write-Output $pr.Text
